# feat: add 2022-Q1 data
#
# The workbook currently has 3 sheets: "2021-Q3", "2021-Q4", "总计".
# This script:
#   1. Inserts a new sheet "2022-Q1" (fund holdings detail) right after
#      "2021-Q4" and before "总计", populated like the other quarterly
#      detail sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#      持有市值(亿元)/仓位排名).
#   2. Adds a new summary row for "2022-Q1" at the top of the "总计"
#      (summary) sheet's data, pushing the existing 2021-Q4 / 2021-Q3
#      rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the workbook's existing bold / centered / thin-bordered
# "header" look (used for column headers and the row-index column) to a
# range, reusing the formatting already used elsewhere in the workbook.
# ---------------------------------------------------------------------
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# =======================================================================
# Step 1: insert the new "2022-Q1" worksheet after "2021-Q4"
# =======================================================================
$sheetQ4 = $wb.Worksheets.Item("2021-Q4")
$sheetQ1 = $wb.Worksheets.Add($null, $sheetQ4)
$sheetQ1.Name = "2022-Q1"

# -- header row --------------------------------------------------------
$sheetQ1.Cells.Item(1, 2).Value = "基金代码"
$sheetQ1.Cells.Item(1, 3).Value = "基金名称"
$sheetQ1.Cells.Item(1, 4).Value = "基金规模"
$sheetQ1.Cells.Item(1, 5).Value = "股票总仓位"
$sheetQ1.Cells.Item(1, 6).Value = "仓位占比"
$sheetQ1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$sheetQ1.Cells.Item(1, 8).Value = "仓位排名"
Set-HeaderStyle $sheetQ1.Range("B1:H1")

# -- row 2: 银河嘉谊灵活配置混合A ---------------------------------------
$sheetQ1.Cells.Item(2, 1).Value = 0
$sheetQ1.Cells.Item(2, 2).Value = "'005459"
$sheetQ1.Cells.Item(2, 3).Value = "银河嘉谊灵活配置混合A"
$sheetQ1.Cells.Item(2, 4).Value = "'6.47"
$sheetQ1.Cells.Item(2, 5).Value = "'39.69"
$sheetQ1.Cells.Item(2, 6).Value = "'0.74"
$sheetQ1.Cells.Item(2, 7).Value = "'0.0479"
$sheetQ1.Cells.Item(2, 8).Value = 4
Set-HeaderStyle $sheetQ1.Range("A2")

# -- row 3: 银河嘉谊灵活配置混合C ---------------------------------------
$sheetQ1.Cells.Item(3, 1).Value = 1
$sheetQ1.Cells.Item(3, 2).Value = "'005460"
$sheetQ1.Cells.Item(3, 3).Value = "银河嘉谊灵活配置混合C"
$sheetQ1.Cells.Item(3, 4).Value = "'2.79"
$sheetQ1.Cells.Item(3, 5).Value = "'39.69"
$sheetQ1.Cells.Item(3, 6).Value = "'0.74"
$sheetQ1.Cells.Item(3, 7).Value = "'0.0206"
$sheetQ1.Cells.Item(3, 8).Value = 4
Set-HeaderStyle $sheetQ1.Range("A3")

# =======================================================================
# Step 2: add the 2022-Q1 row to the "总计" (summary) sheet
# =======================================================================
$sheetTotal = $wb.Worksheets.Item("总计")

# Insert a fresh row right under the header, pushing the 2021-Q4 /
# 2021-Q3 rows down by one.
$sheetTotal.Rows.Item(2).Insert()
$sheetTotal.Range("A2:D2").ClearFormats()

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 2
$sheetTotal.Cells.Item(2, 4).Value = 0.07000000000000001
Set-HeaderStyle $sheetTotal.Range("A2")

# Restore the A-column sequential index (0,1,2) for the rows that got
# shifted down by the insert.
$sheetTotal.Cells.Item(3, 1).Value = 1
$sheetTotal.Cells.Item(4, 1).Value = 2

Write-Host "2022-Q1 data added"
